# Generate Report for Handoff
# Adds two new handed-off files (17f14853-... and 2792364d-...) to the
# localization status workbook, on all three sheets (Overview, zh-cn, de-de),
# ahead of the trailing ".localization-config" row.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/fd616bcc6d04dd1d8c863815ed0990ab2a9a22c2"
$zhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f501402cfda97aa9e55c73e81d57adb22c7b8cca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d76785eeaef48f1713fafdf1a6ef89956886c32d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$file17 = "17f14853-8e61-4ea6-914d-10a1a3ff1c0e.md"
$file27 = "2792364d-b5bc-4750-8fb6-f883abaa8b65.md"
$cfgFile = ".localization-config"

$xlf17zh = "17f14853-8e61-4ea6-914d-10a1a3ff1c0e.2ee94b7c2a9410e8803672872bcb050f9b5ddd00.zh-cn.xlf"
$xlf27zh = "2792364d-b5bc-4750-8fb6-f883abaa8b65.fb48680b57e07b32bf4e2d360af1153df7d062dd.zh-cn.xlf"
$xlf17de = "17f14853-8e61-4ea6-914d-10a1a3ff1c0e.2ee94b7c2a9410e8803672872bcb050f9b5ddd00.de-de.xlf"
$xlf27de = "2792364d-b5bc-4750-8fb6-f883abaa8b65.fb48680b57e07b32bf4e2d360af1153df7d062dd.de-de.xlf"

$readyStatus = "Ready for handoff"
$notLocalized = "Not to be localized"
$inTranslation = "In Translation"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = $file17
$ws.Range("B4").Value = $readyStatus
$ws.Range("C4").Value = $readyStatus

$ws.Range("A5").Value = $file27
$ws.Range("B5").Value = $readyStatus
$ws.Range("C5").Value = $readyStatus

$ws.Range("A6").Value = $cfgFile
$ws.Range("B6").Value = $notLocalized
$ws.Range("C6").Value = $notLocalized

# Hyperlinks are rebuilt from scratch (any delete on this worksheet clears
# all of them), then re-added in row order so relationship ids come out
# sequential again.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/30725227-8aa7-41ee-b39d-a961f8805581.md", "", "", "30725227-8aa7-41ee-b39d-a961f8805581.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md", "", "", "a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file17", "", "", $file17) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file27", "", "", $file27) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/$cfgFile", "", "", $cfgFile) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A4").Value = $file17
$ws.Range("B4").Value = $readyStatus
$ws.Range("C4").Value = $xlf17zh
$ws.Range("D4").Value = "2016-03-07 04:17:44"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = $file27
$ws.Range("B5").Value = $readyStatus
$ws.Range("C5").Value = $xlf27zh
$ws.Range("D5").Value = "2016-03-07 04:17:44"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

$ws.Range("A6").Value = $cfgFile
$ws.Range("B6").Value = $notLocalized
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/30725227-8aa7-41ee-b39d-a961f8805581.md", "", "", "30725227-8aa7-41ee-b39d-a961f8805581.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$zhBase/30725227-8aa7-41ee-b39d-a961f8805581.21ed816ca6caeb7cd93d77316309b4cabacf75c7.zh-cn.xlf", "", "", "30725227-8aa7-41ee-b39d-a961f8805581.21ed816ca6caeb7cd93d77316309b4cabacf75c7.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md", "", "", "a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$zhBase/a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.0f65f2a05226eaebb25f15ca451a906d036a076d.zh-cn.xlf", "", "", "a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.0f65f2a05226eaebb25f15ca451a906d036a076d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file17", "", "", $file17) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "$zhBase/$xlf17zh", "", "", $xlf17zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file27", "", "", $file27) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "$zhBase/$xlf27zh", "", "", $xlf27zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/$cfgFile", "", "", $cfgFile) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A4").Value = $file17
$ws.Range("B4").Value = $readyStatus
$ws.Range("C4").Value = $xlf17de
$ws.Range("D4").Value = "2016-03-07 04:17:55"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = $file27
$ws.Range("B5").Value = $readyStatus
$ws.Range("C5").Value = $xlf27de
$ws.Range("D5").Value = "2016-03-07 04:17:55"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

$ws.Range("A6").Value = $cfgFile
$ws.Range("B6").Value = $notLocalized
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/30725227-8aa7-41ee-b39d-a961f8805581.md", "", "", "30725227-8aa7-41ee-b39d-a961f8805581.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$deBase/30725227-8aa7-41ee-b39d-a961f8805581.21ed816ca6caeb7cd93d77316309b4cabacf75c7.de-de.xlf", "", "", "30725227-8aa7-41ee-b39d-a961f8805581.21ed816ca6caeb7cd93d77316309b4cabacf75c7.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md", "", "", "a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$deBase/a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.0f65f2a05226eaebb25f15ca451a906d036a076d.de-de.xlf", "", "", "a888ceb9-24b5-4a1a-89bc-9ea1bd510dd7.0f65f2a05226eaebb25f15ca451a906d036a076d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file17", "", "", $file17) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "$deBase/$xlf17de", "", "", $xlf17de) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file27", "", "", $file27) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "$deBase/$xlf27de", "", "", $xlf27de) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/$cfgFile", "", "", $cfgFile) | Out-Null

Write-Output "Report generated for handoff: $file17, $file27"
